# Apply the LinuxForHealth rebrand edit to the StructureDefinition workbook.

$wb = $excel.ActiveWorkbook

# ---- Sheet "Metadata" ----
$meta = $wb.Worksheets.Item("Metadata")

# URL
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-termination-date"
# Version
$meta.Range("B3").Value = "8.0.0"
# Date
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
# Publisher
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---- Sheet "Elements" ----
$elements = $wb.Worksheets.Item("Elements")

# Row 5 = Extension.url -> Fixed Value (column Q) contains the same URL
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-termination-date"

# Row 2 = Extension (root element) -> Constraint(s) (column AI) cleared
$elements.Range("AI2").Value = ""
